$wb = $excel.ActiveWorkbook

# Add the missing "Hard" score row to the Hangman sheet.
$ws = $wb.Worksheets.Item("Hangman")
$ws.Cells.Item(4, 1).Value = "Hard"
$ws.Cells.Item(4, 2).Value = 457
$ws.Cells.Item(4, 3).Value = 152
